{"js": "// Load all paragraphs so we can find the exact anchor points for the new\n// content (matching paragraphs by their current text rather than a fixed\n// index, so the script is robust to minor load-order differences).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Hunk 1: a new numbered (\"ListParagraph\" / numId 4) bullet is inserted\n// right after the\n//   Console.WriteLine($\"<text>{variable}\"); - Shorcut cw+tab+tab\n// paragraph (and therefore right before \"<variable> = Console.Readline();\").\n// ---------------------------------------------------------------------\nlet shortcutPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Shorcut\") !== -1 && p.text.indexOf(\"cw+tab+tab\") !== -1) {\n    shortcutPara = p;\n    break;\n  }\n}\n\nif (shortcutPara) {\n  // insertParagraph after this list item automatically inherits its\n  // pStyle=\"ListParagraph\" + numPr (ilvl 0 / numId 4) + rPr(lang en-AU),\n  // exactly matching the paragraph already used by the list above it.\n  shortcutPara.insertParagraph(\n    \"Console.WriteLine(\\u201C<text>\\u201D+<variable>)\",\n    Word.InsertLocation.after\n  );\n}\n\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Hunk 2: two new paragraphs are appended at the very end of the body,\n// after \"int <name> = new int[]{2,4,6};\":\n//   1) \"<name>.arraylenth \u2013 length of array\"   (Consolas / 000000 / 9.5pt)\n//   2) an empty paragraph (same formatting) that now carries the\n//      \"_GoBack\" bookmark, which used to sit at the end of the old last\n//      paragraph.\n// ---------------------------------------------------------------------\nparagraphs.load(\"items,text\");\nawait context.sync();\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// The \"_GoBack\" bookmark currently lives at the end of the last paragraph;\n// remove it there so it can be re-created at the new end of the document.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nconst arrayLenParagraph = lastParagraph.insertParagraph(\n  \"<name>.arraylenth \\u2013 length of array\",\n  Word.InsertLocation.after\n);\narrayLenParagraph.font.set({ name: \"Consolas\", color: \"#000000\", size: 9.5 });\n\nconst trailingParagraph = arrayLenParagraph.insertParagraph(\"\", Word.InsertLocation.after);\ntrailingParagraph.font.set({ name: \"Consolas\", color: \"#000000\", size: 9.5 });\n\ntrailingParagraph.getRange(\"Start\").insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Hunk 1: a new numbered (\"ListParagraph\" / numId 4) bullet is inserted\n# right after the\n#   Console.WriteLine($\"<text>{variable}\"); - Shorcut cw+tab+tab\n# paragraph (and therefore right before \"<variable> = Console.Readline();\").\n# ---------------------------------------------------------------------\n$findRange = $d.Content\n$findRange.Find.Execute(\"cw+tab+tab\") | Out-Null\n$shortcutPara = $findRange.Paragraphs(1)\n\n# InsertParagraphAfter inherits the paragraph's pStyle=\"ListParagraph\" +\n# numPr (ilvl 0 / numId 4) + rPr(lang en-AU), exactly matching the list\n# item already used above it.\n$shortcutPara.Range.InsertParagraphAfter()\n$newListPara = $shortcutPara.Next()\n\n$leftQuote = [char]8220\n$rightQuote = [char]8221\n$newListPara.Range.Text = \"Console.WriteLine($leftQuote<text>$rightQuote+<variable>)\"\n\n# ---------------------------------------------------------------------\n# Hunk 2: two new paragraphs are appended at the very end of the body,\n# after \"int <name> = new int[]{2,4,6};\":\n#   1) \"<name>.arraylenth - length of array\"   (Consolas / 000000 / 9.5pt)\n#   2) an empty paragraph (same formatting) that now carries the\n#      \"_GoBack\" bookmark, which used to sit at the end of the old last\n#      paragraph.\n# ---------------------------------------------------------------------\n$d.Bookmarks(\"_GoBack\").Delete()\n\n$count = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs($count)\n$lastPara.Range.InsertParagraphAfter()\n\n$arrayLenPara = $d.Paragraphs($count + 1)\n$arrayLenPara.Range.Font.Name = \"Consolas\"\n$arrayLenPara.Range.Font.Color = 0\n$arrayLenPara.Range.Font.Size = 9.5\n$enDash = [char]8211\n$arrayLenPara.Range.Text = \"<name>.arraylenth $enDash length of array\"\n\n$arrayLenPara.Range.InsertParagraphAfter()\n$trailingPara = $d.Paragraphs($count + 2)\n$trailingPara.Range.Font.Name = \"Consolas\"\n$trailingPara.Range.Font.Color = 0\n$trailingPara.Range.Font.Size = 9.5\n\n$bookmarkRange = $trailingPara.Range.Duplicate\n$bookmarkRange.Collapse(1)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n"}
